$wb = $excel.ActiveWorkbook

# --- Step 1: duplicate "trial 5" (in its original, unedited state) to create "trial 6" ---
$ws5 = $wb.Worksheets.Item("trial 5")
$ws5.Copy($null, $ws5) | Out-Null
$ws6 = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws6.Name = "trial 6"

# --- Step 2: correct a couple of mislabeled microsite names and a typo on "trial 5" ---
$ws5.Range("A3").Value = "open 3"
$ws5.Range("A5").Value = "open 2"
$ws5.Range("C7").Value = 15491

# --- Step 3: update data on the new "trial 6" sheet ---
$ws6.Range("A2").Value = "burlap 1"
$ws6.Range("B2").Value = 20881725
$ws6.Range("C2").Value = 15406

$ws6.Range("A3").Value = "open1"
$ws6.Range("B3").Value = 20881726
$ws6.Range("C3").Value = 15512

$ws6.Range("A4").Value = "canvas 1"
$ws6.Range("B4").Value = "L1"
$ws6.Range("C4").Value = 15289

$ws6.Range("A5").Value = "open 2"
$ws6.Range("B5").Value = 20884543
$ws6.Range("C5").Value = 15491

$ws6.Range("A6").Value = "canvas 2"
$ws6.Range("B6").Value = 20881725
$ws6.Range("C6").Value = 15402

$ws6.Range("A7").Value = "open 3"
$ws6.Range("B7").Value = 20881733
$ws6.Range("C7").Value = 15531

# --- Step 4: the original "trial 5" sheet is no longer the selected tab; ---
# --- whole first row is left selected there, "trial 6" becomes the active tab ---
$ws5.Activate() | Out-Null
$ws5.Rows(1).Select() | Out-Null

$ws6.Activate() | Out-Null
$ws6.Range("D7").Select() | Out-Null

Write-Output "done"
